$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns F ("Type *") and G ("Folio No") entirely.
# This shifts old column H ("Tag") -> F and old column I ("Instrument *") -> G,
# carrying their comments and formatting along automatically.
$ws.Range("F1:G1").EntireColumn.Delete()

# Leave selection on column F, matching the post-delete cursor position.
$ws.Range("F1:F1048576").Select()
